# 氧化铝.xlsx edit
#
# Two structural changes, derived from the OOXML diff:
#
# 1) The data is organised in 4-row "year" blocks (rows 2-69), each block
#    being quarter labels A/B/C/D (e.g. "2000年A".."2000年D"). Within every
#    block the two middle rows (the "B" and "C" quarters) have had their
#    entire row contents (columns A-E) swapped with each other, while the
#    "A" and "D" rows (first/last of each block) stay in place.
#
# 2) Columns F ("氧化铝产销率") and G ("氧化铝销售量") are deleted entirely
#    (duplicate/derived data also present in columns B and E), shrinking the
#    sheet's dimension from A1:G69 to A1:E69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the "B" and "C" quarter rows inside each year block ---
for ($blockStart = 2; $blockStart -le 69; $blockStart += 4) {
    $rowB = $blockStart + 1   # the "...年B" row of this block
    $rowC = $blockStart + 2   # the "...年C" row of this block

    $bValues = @()
    $cValues = @()
    for ($col = 1; $col -le 5; $col++) {
        $bValues += ,$ws.Cells.Item($rowB, $col).Value2
        $cValues += ,$ws.Cells.Item($rowC, $col).Value2
    }
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($rowB, $col).Value2 = $cValues[$col - 1]
        $ws.Cells.Item($rowC, $col).Value2 = $bValues[$col - 1]
    }
}

# --- Step 2: delete columns F and G (氧化铝产销率, 氧化铝销售量) ---
$ws.Range("F:G").Delete()
